$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("N2").Value = 7152.019986098921
$ws1.Range("O2").Value = 6979.505869462281

$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = 5707.815717280662
$ws2.Range("I2").Value = 44492.05901988943
$ws2.Range("L2").Value = 66334.06707325629
$ws2.Range("M2").Value = 21991.42050229464
$ws2.Range("O2").Value = 12076.80007217423

$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("A2").Value = 2927.360317916481
$ws3.Range("B2").Value = 7940.887964949257
$ws3.Range("E2").Value = 67179.99183625776
$ws3.Range("I2").Value = 59530.75343380851
$ws3.Range("L2").Value = 66334.06707325629
$ws3.Range("M2").Value = 25547.11936466757
$ws3.Range("N2").Value = 15110.18723133435
$ws3.Range("O2").Value = 14758.23231153656

$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("A2").Value = 2927.360317916481
$ws4.Range("B2").Value = 7940.887964949257
$ws4.Range("E2").Value = 67179.99183625776
$ws4.Range("I2").Value = 59530.75343380851
$ws4.Range("L2").Value = 66334.06707325629
$ws4.Range("M2").Value = 25547.11936466757
$ws4.Range("N2").Value = 15217.15884705062
$ws4.Range("O2").Value = 14758.23231153656

$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("A2").Value = 6352.985609279765
$ws5.Range("B2").Value = 7940.887964949257
$ws5.Range("E2").Value = 67179.99183625776
$ws5.Range("I2").Value = 59530.75343380851
$ws5.Range("L2").Value = 66334.06707325629
$ws5.Range("M2").Value = 25547.11936466757
$ws5.Range("N2").Value = 15760.98179418098
$ws5.Range("O2").Value = 17091.72736868646

$ws6 = $wb.Worksheets.Item(6)
$ws6.Range("A2").Value = 6352.985609279765
$ws6.Range("B2").Value = 7940.887964949257
$ws6.Range("E2").Value = 67179.99183625776
$ws6.Range("I2").Value = 59530.75343380851
$ws6.Range("L2").Value = 66334.06707325629
$ws6.Range("M2").Value = 25547.11936466757
$ws6.Range("N2").Value = 15760.98179418098
$ws6.Range("O2").Value = 17091.72736868646
